$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

$data = New-Object 'object[,]' 24,14
$data[0,0] = 0.1423503261597858
$data[0,1] = 0
$data[0,2] = 0.234267617884484
$data[0,3] = 0.1793056353713887
$data[0,4] = 1.123359372945323
$data[0,5] = 0.5551713029926262
$data[0,6] = 0.6982842392986512
$data[0,7] = 0
$data[0,8] = 0.185453310251404
$data[0,9] = 0.6090961603234462
$data[0,10] = 0
$data[0,11] = 0.2392297715950988
$data[0,12] = 0
$data[0,13] = 2.480822137930886
$data[1,0] = 0.1329786013034067
$data[1,1] = 0
$data[1,2] = 0.2311197359476722
$data[1,3] = 0.1789700496125164
$data[1,4] = 1.129610054076352
$data[1,5] = 0.5597072547692648
$data[1,6] = 0.7043583452538513
$data[1,7] = 0
$data[1,8] = 0.1865001412834566
$data[1,9] = 0.5316092208624923
$data[1,10] = 0
$data[1,11] = 0.2193315339418049
$data[1,12] = 0
$data[1,13] = 2.502893528269482
$data[2,0] = 0.1272939203481798
$data[2,1] = 0
$data[2,2] = 0.2292743115819462
$data[2,3] = 0.1788473459038009
$data[2,4] = 1.134098854262731
$data[2,5] = 0.5629032805972258
$data[2,6] = 0.7084115039277634
$data[2,7] = 0
$data[2,8] = 0.1872410927430508
$data[2,9] = 0.4838346948781691
$data[2,10] = 0
$data[2,11] = 0.207128646995109
$data[2,12] = 0
$data[2,13] = 2.517985059048812
$data[3,0] = 0.124995031013043
$data[3,1] = 0
$data[3,2] = 0.2285443640525813
$data[3,3] = 0.17881833937086
$data[3,4] = 1.136091789504363
$data[3,5] = 0.5643089328828879
$data[3,6] = 0.7101446230665331
$data[3,7] = 0
$data[3,8] = 0.1875677383262193
$data[3,9] = 0.4643179745582415
$data[3,10] = 0
$data[3,11] = 0.2021599195801116
$data[3,12] = 0
$data[3,13] = 2.524522029385423
$data[4,0] = 0.1246143738932943
$data[4,1] = 0
$data[4,2] = 0.2284244933158561
$data[4,3] = 0.178814792091373
$data[4,4] = 1.136432604177813
$data[4,5] = 0.5645485736162641
$data[4,6] = 0.7104373252839409
$data[4,7] = 0
$data[4,8] = 0.1876234698130652
$data[4,9] = 0.4610743668931434
$data[4,10] = 0
$data[4,11] = 0.2013351209623693
$data[4,12] = 0
$data[4,13] = 2.525630862252626
$data[5,0] = 0.127262844984088
$data[5,1] = 0
$data[5,2] = 0.2292643777238226
$data[5,3] = 0.1788468696482717
$data[5,4] = 1.134125068752695
$data[5,5] = 0.5629218197928623
$data[5,6] = 0.7084345476377507
$data[5,7] = 0
$data[5,8] = 0.1872453979686064
$data[5,9] = 0.4835716789375226
$data[5,10] = 0
$data[5,11] = 0.207061620145609
$data[5,12] = 0
$data[5,13] = 2.518071651901934
$data[6,0] = 0.139104605869548
$data[6,1] = 0
$data[6,2] = 0.2331641514853828
$data[6,3] = 0.1791726488709422
$data[6,4] = 1.125379550974309
$data[6,5] = 0.556649957561838
$data[6,6] = 0.700311447447092
$data[6,7] = 0
$data[6,8] = 0.1857938837283051
$data[6,9] = 0.5824204031498539
$data[6,10] = 0
$data[6,11] = 0.2323660085554948
$data[6,12] = 0
$data[6,13] = 2.488112748237654
$data[7,0] = 0.1628725636555117
$data[7,1] = 0
$data[7,2] = 0.2415010115346661
$data[7,3] = 0.1804713915996778
$data[7,4] = 1.11339245824847
$data[7,5] = 0.5476156562023817
$data[7,6] = 0.6869480570104614
$data[7,7] = 0
$data[7,8] = 0.1837262308379479
$data[7,9] = 0.7746459783755313
$data[7,10] = 0
$data[7,11] = 0.282091728160772
$data[7,12] = 0
$data[7,13] = 2.441584801523604
$data[8,0] = 0.1806618691844619
$data[8,1] = 0
$data[8,2] = 0.2480419082616834
$data[8,3] = 0.1818262799404735
$data[8,4] = 1.107732175530067
$data[8,5] = 0.5429747890700156
$data[8,6] = 0.6786917608833818
$data[8,7] = 0
$data[8,8] = 0.1826815563245248
$data[8,9] = 0.9148329650722644
$data[8,10] = 0
$data[8,11] = 0.3186749555087118
$data[8,12] = 0
$data[8,13] = 2.414859744543492
$data[9,0] = 0.1888245138290188
$data[9,1] = 0
$data[9,2] = 0.2511069222266684
$data[9,3] = 0.1825293472030793
$data[9,4] = 1.10584033831978
$data[9,5] = 0.5412984688520481
$data[9,6] = 0.6752743678819684
$data[9,7] = 0
$data[9,8] = 0.182309266818578
$data[9,9] = 0.9783706907473686
$data[9,10] = 0
$data[9,11] = 0.3353257674520762
$data[9,12] = 0
$data[9,13] = 2.404323406680902
$data[10,0] = 0.1919254478107177
$data[10,1] = 0
$data[10,2] = 0.2522803438699981
$data[10,3] = 0.1828080177991502
$data[10,4] = 1.105222145185003
$data[10,5] = 0.5407263224634136
$data[10,6] = 0.6740289187424509
$data[10,7] = 0
$data[10,8] = 0.1821830857952449
$data[10,9] = 1.002395914907311
$data[10,10] = 0
$data[10,11] = 0.3416319765200058
$data[10,12] = 0
$data[10,13] = 2.400566817593017
$data[11,0] = 0.1912571679307149
$data[11,1] = 0
$data[11,2] = 0.252027060388528
$data[11,3] = 0.1827474486601126
$data[11,4] = 1.105350916760088
$data[11,5] = 0.5408467567522734
$data[11,6] = 0.6742949854433533
$data[11,7] = 0
$data[11,8] = 0.1822096031267506
$data[11,9] = 0.9972232390592808
$data[11,10] = 0
$data[11,11] = 0.3402737871324177
$data[11,12] = 0
$data[11,13] = 2.401365487454143
$data[12,0] = 0.1890794317851032
$data[12,1] = 0
$data[12,2] = 0.2512032050649537
$data[12,3] = 0.1825520246068635
$data[12,4] = 1.105787511139127
$data[12,5] = 0.5412501419996829
$data[12,6] = 0.675170929279858
$data[12,7] = 0
$data[12,8] = 0.1822985893087719
$data[12,9] = 0.9803479729719413
$data[12,10] = 0
$data[12,11] = 0.3358445670112289
$data[12,12] = 0
$data[12,13] = 2.404009673199141
$data[13,0] = 0.1877467918203735
$data[13,1] = 0
$data[13,2] = 0.2507002294184986
$data[13,3] = 0.1824339399656871
$data[13,4] = 1.106067726386243
$data[13,5] = 0.5415053879247722
$data[13,6] = 0.6757138042847259
$data[13,7] = 0
$data[13,8] = 0.1823550227480837
$data[13,9] = 0.9700067660369029
$data[13,10] = 0
$data[13,11] = 0.3331316489456739
$data[13,12] = 0
$data[13,13] = 2.405659700653672
$data[14,0] = 0.1801298186364733
$data[14,1] = 0
$data[14,2] = 0.2478433948002987
$data[14,3] = 0.1817820748132597
$data[14,4] = 1.107869554842445
$data[14,5] = 0.5430931000343548
$data[14,6] = 0.6789219037732579
$data[14,7] = 0
$data[14,8] = 0.1827079572353441
$data[14,9] = 0.9106757891918278
$data[14,10] = 0
$data[14,11] = 0.3175869332068828
$data[14,12] = 0
$data[14,13] = 2.41558095464157
$data[15,0] = 0.1754749063575218
$data[15,1] = 0
$data[15,2] = 0.2461136735614389
$data[15,3] = 0.1814043624610591
$data[15,4] = 1.109149852488457
$data[15,5] = 0.5441785608062588
$data[15,6] = 0.6809766341549945
$data[15,7] = 0
$data[15,8] = 0.182950832213578
$data[15,9] = 0.8742171587700227
$data[15,10] = 0
$data[15,11] = 0.3080527709279721
$data[15,12] = 0
$data[15,13] = 2.422082643440092
$data[16,0] = 0.1728041421246473
$data[16,1] = 0
$data[16,2] = 0.2451272167508165
$data[16,3] = 0.1811952771417715
$data[16,4] = 1.109950539229516
$data[16,5] = 0.5448438111600069
$data[16,6] = 0.6821903174667199
$data[16,7] = 0
$data[16,8] = 0.183100216840522
$data[16,9] = 0.8532251779980982
$data[16,10] = 0
$data[16,11] = 0.3025698346939834
$data[16,12] = 0
$data[16,13] = 2.425974810076525
$data[17,0] = 0.1719010093194129
$data[17,1] = 0
$data[17,2] = 0.2447946710624649
$data[17,3] = 0.1811258878644253
$data[17,4] = 1.110232680975976
$data[17,5] = 0.5450760784884281
$data[17,6] = 0.6826067219251115
$data[17,7] = 0
$data[17,8] = 0.1831524602297776
$data[17,9] = 0.846113930850521
$data[17,10] = 0
$data[17,11] = 0.3007135661314564
$data[17,12] = 0
$data[17,13] = 2.427318828493597
$data[18,0] = 0.175969746040451
$data[18,1] = 0
$data[18,2] = 0.246296933319158
$data[18,3] = 0.1814437258180739
$data[18,4] = 1.109006908617076
$data[18,5] = 0.544058775567855
$data[18,6] = 0.6807546077146824
$data[18,7] = 0
$data[18,8] = 0.1829239749841207
$data[18,9] = 0.8781005237553074
$data[18,10] = 0
$data[18,11] = 0.3090676123088301
$data[18,12] = 0
$data[18,13] = 2.421374736033371
$data[19,0] = 0.1897188179324871
$data[19,1] = 0
$data[19,2] = 0.2514448456897895
$data[19,3] = 0.1826090882147255
$data[19,4] = 1.105656607691351
$data[19,5] = 0.5411299570715471
$data[19,6] = 0.6749123233604024
$data[19,7] = 0
$data[19,8] = 0.1822720503698534
$data[19,9] = 0.9853056129479398
$data[19,10] = 0
$data[19,11] = 0.3371455146632911
$data[19,12] = 0
$data[19,13] = 2.403226678718198
$data[20,0] = 0.1987623562361449
$data[20,1] = 0
$data[20,2] = 0.2548836726137438
$data[20,3] = 0.183443179532663
$data[20,4] = 1.104039394524804
$data[20,5] = 0.5395809676030723
$data[20,6] = 0.6713775730178355
$data[20,7] = 0
$data[20,8] = 0.1819322251639406
$data[20,9] = 1.055164869906321
$data[20,10] = 0
$data[20,11] = 0.3555011663437355
$data[20,12] = 0
$data[20,13] = 2.392725773675011
$data[21,0] = 0.1939304263941324
$data[21,1] = 0
$data[21,2] = 0.2530415364294498
$data[21,3] = 0.1829913912317274
$data[21,4] = 1.10485016299279
$data[21,5] = 0.5403742436381833
$data[21,6] = 0.6732382014347422
$data[21,7] = 0
$data[21,8] = 0.1821057066668175
$data[21,9] = 1.017898956497163
$data[21,10] = 0
$data[21,11] = 0.3457040645729776
$data[21,12] = 0
$data[21,13] = 2.398205810720867
$data[22,0] = 0.1757460121652912
$data[22,1] = 0
$data[22,2] = 0.2462140567090216
$data[22,3] = 0.1814259045154607
$data[22,4] = 1.109071332233341
$data[22,5] = 0.544112802133732
$data[22,6] = 0.6808548849686673
$data[22,7] = 0
$data[22,8] = 0.1829360867589358
$data[22,9] = 0.8763449522310509
$data[22,10] = 0
$data[22,11] = 0.3086088075368494
$data[22,12] = 0
$data[22,13] = 2.421694300479572
$data[23,0] = 0.1563848411108353
$data[23,1] = 0
$data[23,2] = 0.2391723266877932
$data[23,3] = 0.1800495649141745
$data[23,4] = 1.116082642526145
$data[23,5] = 0.5497094966107028
$data[23,6] = 0.6902887613301445
$data[23,7] = 0
$data[23,8] = 0.1842022469703224
$data[23,9] = 0.722822712579557
$data[23,10] = 0
$data[23,11] = 0.2686298917504004
$data[23,12] = 0
$data[23,13] = 2.452862569534048

$ws.Range("B2:O25").Value = $data
